# "SP and process file"
# Update the branch-code (C) column from the text "195" to the numeric value 19,
# and populate the Classification_TYPE (J) column with numeric flags for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C: Brn_Code - change from text "195" to numeric 19
$ws.Range("C2").Value = 19
$ws.Range("C3").Value = 19
$ws.Range("C4").Value = 19

# Column J: Classification_TYPE - fill in numeric classification values
$ws.Range("J2").Value = 1
$ws.Range("J3").Value = 0
$ws.Range("J4").Value = 0

# Restore the active cell selection as last left by the author
$ws.Range("J14").Select() | Out-Null
